# kane_census_mock.xlsx — bring header row in line with the "current format
# of provided file": every demographic header (all columns except Date and
# GH-Male) gains a single trailing space, and the sheet's saved view/
# selection moves from A10/D20 to B1(ish)/U1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): append a trailing space to every header except
#     A1 ("Date") and B1 ("GH-Male"), which are left untouched. Rewriting
#     each cell's Value causes the shared-string table to be rebuilt with
#     the new (space-padded) text, matching the target sharedStrings.xml.
$ws.Range("C1").Value = "GH-Female "
$ws.Range("D1").Value = "GH-White "
$ws.Range("E1").Value = "GH - Black "
$ws.Range("F1").Value = "GH-Other "
$ws.Range("G1").Value = "Scott-Male "
$ws.Range("H1").Value = "Scott-Female "
$ws.Range("I1").Value = "Scott-White "
$ws.Range("J1").Value = "Scott-Black "
$ws.Range("K1").Value = "Scott-Other "
$ws.Range("L1").Value = "McK-Male "
$ws.Range("M1").Value = "McK-Female "
$ws.Range("N1").Value = "McK-White "
$ws.Range("O1").Value = "McK-Black "
$ws.Range("P1").Value = "McK-Other "
$ws.Range("Q1").Value = "Ross-Male "
$ws.Range("R1").Value = "Ross-Female "
$ws.Range("S1").Value = "Ross-White "
$ws.Range("T1").Value = "Ross-Black "
$ws.Range("U1").Value = "Ross-Other "

# --- View/selection: the saved sheet view moves its selection to U1
#     (top-right header cell), matching the diff's
#     <selection activeCell="U1" sqref="U1"/>.
$ws.Activate()
$ws.Range("U1").Select()
